$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("U2").Value = 8130.35
$ws.Range("V2").Value = 1208
$ws.Range("AG2").Value = 227433.75

# Row 3 - Bibi Cell Vieiralves
$ws.Range("U3").Value = 4136
$ws.Range("V3").Value = 5429
$ws.Range("AG3").Value = 106284.2

# Row 4 - Bibi Cell Manauara
$ws.Range("U4").Value = 3561.5
$ws.Range("V4").Value = 4526
$ws.Range("W4").Value = 2936
$ws.Range("AG4").Value = 69753.60000000001

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("U5").Value = 1795
$ws.Range("V5").Value = 2338.01
$ws.Range("W5").Value = 1118.5
$ws.Range("AG5").Value = 59630.7

# Row 6 - total
$ws.Range("U6").Value = 17622.85
$ws.Range("V6").Value = 13501.01
$ws.Range("W6").Value = 4054.5
$ws.Range("AG6").Value = 463102.25
